$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Pour la prochaine fois")

# Clear out the previous content in the used range (B:C) before rewriting
$ws.Range("B2:C16").ClearContents()

# New task breakdown (skin animation work split up, new items added)
# Write order matters for shared-string allocation order, so follow the
# exact sequence the new strings were introduced in.
$ws.Range("B2").Value = "joueur enemie haut dessus de la box "
$ws.Range("C2").Value = "2.8 environ de trop en z, psq skin est sur la box de detection "

$ws.Range("B4").Value = "Animation perso (marcher, sauter, mort…)"

$ws.Range("B5").Value = "Entity enemi : RegisterHandler_ for hit, death…"
$ws.Range("B6").Value = "deleteRemotePlayer : enlever le perso de la scene"
$ws.Range("C5").Value = "pour les différentes animatiosn. Mais la fonction qui envoie des msg entre les composant a pas l'air de marcher de fou"

$ws.Range("C14").Value = "pour une prochiane co, sauvegarder dans le temps"

$ws.Range("B8").Value = "afficher pseudo sur skins"

$ws.Range("B10").Value = "POV d'un joueur bug desfois (mouvement souris trop brusque ?)"

$ws.Range("B12").Value = "optimiser le code en general"

$ws.Range("B14").Value = "Sauvegarder kd par joueurs ? "

$ws.Range("B16").Value = "Faille xss/autre securités"

$ws.Range("B8").Select()
